$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 5 & 6 (Character Movement / Camera work logged more hours) ---

# C5: was a plain value (0.5); now a formula that still evaluates to 1.5
$ws.Range("C5").Formula = "=0.5+1"

# C6: formula gains an extra addend, evaluates to 2.25
$ws.Range("C6").Formula = "=0.25+1.5+0.5"

# E5:E6 get "Date Completed" entries now (shared formula), matching the
# look/format of the existing Date Completed column (copy format from E4
# first, then fill in the formula so the number format / font / alignment
# match column E rather than picking up a default style).
$ws.Range("E4").Copy()
$ws.Range("E5:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E5:E6").Formula = "=DATE(2021,8,18)"

# D5:D6 flip to the "Solved" (Good/green) look now that they're done
$ws.Range("D5:D6").Style = "Good"

# --- New backlog rows 7-10 ---
# Typed in this order so new shared-string entries land in the same order
# as the authored workbook (Character animation, Inventory, Player gear,
# then Fishing system last).
$ws.Range("A8").Value = "Character animation - bare bones run anim"
$ws.Range("A9").Value = "Inventory system"
$ws.Range("A10").Value = "Player gear system"
$ws.Range("A7").Value = "Fishing system - written in outline"

# New rows aren't solved yet -> "Bad" (red) look in the Solved column
$ws.Range("D7:D10").Style = "Bad"

# Move the active selection to I7 (next empty Legend row)
$ws.Range("I7").Select() | Out-Null
